# Add Figure (from Pandoc 3)
#
# Splits the old "Abstract" paragraph style in two: a new "AbstractTitle"
# style (bold, centered, colored "Abstract" label) used for a heading
# paragraph that precedes the existing "Abstract" body paragraph.

$d = $word.ActiveDocument

# 1. Define the new "AbstractTitle" paragraph style.
$abstractTitleStyle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitleStyle.NameLocal = "Abstract Title"
$abstractTitleStyle.BaseStyle = "Normal"
$abstractTitleStyle.NextParagraphStyle = "Abstract"
$abstractTitleStyle.QuickStyle = $true

$abstractTitleStyle.ParagraphFormat.KeepWithNext = $true
$abstractTitleStyle.ParagraphFormat.KeepTogether = $true
$abstractTitleStyle.ParagraphFormat.SpaceBefore = 15
$abstractTitleStyle.ParagraphFormat.SpaceAfter = 0
$abstractTitleStyle.ParagraphFormat.Alignment = 1

$abstractTitleStyle.Font.Bold = $true
$abstractTitleStyle.Font.Color = 9067060
$abstractTitleStyle.Font.Size = 10
$abstractTitleStyle.Font.SizeBi = 10

# 2. The existing "Abstract" body style now follows a heading, so it
#    needs less space above it (300 twips -> 100 twips = 5pt).
$abstractStyle = $d.Styles.Item("Abstract")
$abstractStyle.ParagraphFormat.SpaceBefore = 5

# 3. Find the existing "Abstract" body paragraph ("Quite a long
#    description ...") and insert a new "Abstract" heading paragraph
#    right before it.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Abstract") {
        $targetPara = $p
        break
    }
}

$targetIndex = $targetPara.Index
$targetPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = "Abstract"
$newPara.Style = "AbstractTitle"

# 4. Template refresh that came along with the same Pandoc 3 reference
#    doc update: a new "Footnote Block Text" style (a block-text variant
#    scoped to footnotes) and restyled Pandoc syntax-highlighting token
#    styles (ImportTok/BuiltInTok now render in the same green as the
#    other type-ish tokens).
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = "FootnoteText"
$footnoteBlockText.NextParagraphStyle = "FootnoteText"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

$importTok = $d.Styles.Item("ImportTok")
$importTok.Font.Bold = $true
$importTok.Font.Color = 32768

$builtInTok = $d.Styles.Item("BuiltInTok")
$builtInTok.Font.Color = 32768
